$d = $word.ActiveDocument

# Locate the first paragraph's range (the one ending with
# "This is a Microsoft word document.") and append the new text as
# additional runs, right after the existing run.
$p = $d.Paragraphs(1)
$r = $p.Range
$r.Collapse(0)  # wdCollapseEnd -> collapses to the end of the paragraph's range (before the paragraph mark)
$r.MoveEnd(1, -1) | Out-Null  # step back over the paragraph mark so we insert before it

$r.InsertAfter(" (")
$r.Collapse(0)
$r.InsertAfter("Changed main")
$r.Collapse(0)
$r.InsertAfter(")")
